$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Accuracy of Financial Insights -> Accuracy of financial insights, score 9 -> 7
$ws.Range("A3").Value = "*   **Accuracy of financial insights"
$ws.Range("B3").Value = "7"

# Row 4: Value to Investors -> Value to investors, score 9 -> 8
$ws.Range("A4").Value = "*   **Value to investors"
$ws.Range("B4").Value = "8"

# Row 5: Clarity of Writing -> Clarity of writing (score unchanged)
$ws.Range("A5").Value = "*   **Clarity of writing"

# Row 6: score 8 -> 9
$ws.Range("B6").Value = "9"

# Row 7: Accuracy of Financial Insights -> Accuracy of financial insights (score unchanged)
$ws.Range("A7").Value = "*   **Accuracy of financial insights"

# Row 8: Value to Investors -> Value to investors (score unchanged)
$ws.Range("A8").Value = "*   **Value to investors"

# Row 9: Clarity of Writing -> Clarity of writing (score unchanged)
$ws.Range("A9").Value = "*   **Clarity of writing"
